# Redis Release Notes - apply "persistence-available flag" update
#
# Summary of changes (see commit message / diff):
#  1. Title paragraph: wrap "MSOpenTech" in a spell-check proofErr span
#     (splitting the leading run) and move the "_GoBack" bookmark from its
#     old location into the title, splitting " Release Notes" into " R" /
#     "elease Notes" around the bookmark.
#  2. Remove the old "_GoBack" bookmark pair that used to sit right after
#     "...2.8.14" in the "What's new" paragraph.
#  3. Replace the single "Since there have been..." paragraph with:
#       - a "Network layer changes" Heading3
#       - a reworded/re-split paragraph about the networking changes
#       - a new "persistence-available flag" Heading3
#       - a paragraph describing the new flag
#       - an indented paragraph showing the flag syntax
#
# We use Range.InsertXML with literal OOXML so run splits, <w:proofErr/>
# markers and the relocated <w:bookmarkStart/End/> land exactly where the
# target markup puts them (none of that is reachable through the plain
# Word object model properties).

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------
# 1. Title paragraph (paragraph 1)
# ---------------------------------------------------------------------
$titleXml = "<w:p $wNs w:rsidR='004A0C03' w:rsidRDefault='005E44B2' w:rsidP='004A0C03'>" +
  "<w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr></w:pPr>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr><w:t>MSOpenTech</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r w:rsidR='00853344'><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr><w:t>Redis 2.8.14</w:t></w:r>" +
  "<w:r w:rsidR='004A0C03'><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr><w:t xml:space='preserve'> R</w:t></w:r>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
  "<w:r><w:rPr><w:rFonts w:eastAsia='Calibri'/></w:rPr><w:t>elease Notes</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(1).Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2. Drop the old bookmark that used to follow "...2.8.14" in the
#    "What's new" paragraph (8th paragraph: "Our last official release...").
# ---------------------------------------------------------------------
$whatsNewXml = "<w:p $wNs w:rsidR='00010072' w:rsidRDefault='00A10041' w:rsidP='00580900'>" +
  "<w:r><w:t>Our</w:t></w:r>" +
  "<w:r w:rsidR='000B3594'><w:t xml:space='preserve'> last official release was 2.8.12</w:t></w:r>" +
  "<w:r><w:t>. We have merged in the changes up to 2.8.</w:t></w:r>" +
  "<w:r w:rsidR='00853344'><w:t>14</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>. Please see the </w:t></w:r>" +
  "<w:hyperlink r:id='rId4' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships' w:history='1'>" +
    "<w:r w:rsidRPr='00DB18A6'><w:rPr><w:rStyle w:val='Hyperlink'/></w:rPr><w:t xml:space='preserve'>release notes for the </w:t></w:r>" +
    "<w:r w:rsidR='00DB18A6'><w:rPr><w:rStyle w:val='Hyperlink'/></w:rPr><w:t>UNIX</w:t></w:r>" +
    "<w:r w:rsidR='00DB18A6' w:rsidRPr='00DB18A6'><w:rPr><w:rStyle w:val='Hyperlink'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r w:rsidRPr='00DB18A6'><w:rPr><w:rStyle w:val='Hyperlink'/></w:rPr><w:t>2.8 branch</w:t></w:r>" +
  "</w:hyperlink>" +
  "<w:r><w:t xml:space='preserve'> to understand how this impacts Redis functionality.</w:t></w:r>" +
  "<w:r w:rsidR='00DB18A6'><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r w:rsidR='000B3594'><w:t xml:space='preserve'> </w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(8).Range.InsertXML($whatsNewXml)

# ---------------------------------------------------------------------
# 3. Replace the "Since there have been..." paragraph (now paragraph 9)
#    with the new Heading3 + body, then append three more new paragraphs
#    after it for the persistence-available flag section.
# ---------------------------------------------------------------------
$networkHeadingXml = "<w:p $wNs w:rsidR='000B3594' w:rsidRDefault='000B3594' w:rsidP='00580900'>" +
  "<w:pPr><w:pStyle w:val='Heading3'/></w:pPr>" +
  "<w:r><w:t>Network layer changes</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(9).Range.InsertXML($networkHeadingXml)

$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$networkBodyXml = "<w:p $wNs>" +
  "<w:r><w:t>T</w:t></w:r>" +
  "<w:r><w:t>here have been significant changes to the ne</w:t></w:r>" +
  "<w:r><w:t>tworking layer for this version. L</w:t></w:r>" +
  "<w:r><w:t>ikely there will be a few weeks before there is another official (</w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/><w:r><w:t>Chocolatey</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t xml:space='preserve'> and </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/><w:r><w:t>Nuget</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t>) release.</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> Most of these changes target IPv6.</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(10).Range.InsertXML($networkBodyXml)

$p10 = $d.Paragraphs(10)
$p10.Range.InsertParagraphAfter()
$flagHeadingXml = "<w:p $wNs>" +
  "<w:pPr><w:pStyle w:val='Heading3'/></w:pPr>" +
  "<w:proofErr w:type='gramStart'/><w:r><w:t>persistence-available</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:t xml:space='preserve'> flag </w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(11).Range.InsertXML($flagHeadingXml)

$p11 = $d.Paragraphs(11)
$p11.Range.InsertParagraphAfter()
$flagBodyXml = "<w:p $wNs>" +
  "<w:r><w:t xml:space='preserve'>If Redis is to be used as an in-memory-only cache without any kind of persistence, then the </w:t></w:r>" +
  "<w:proofErr w:type='gramStart'/><w:r><w:t>fork(</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:t>) mechanism used by the background AOF/RDB persistence is unnecessary. As an optimization, all persistence can be</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:t>turned off in the Windows version of Redis</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> in this scenario</w:t></w:r>" +
  "<w:r><w:t>. This will disable the creation of the memory mapped heap file, redirect heap allocations to the system heap allocator, and disable commands that would otherwise cause fork() operations:</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> B</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>GSAVE and BGREWRITEAOF. This flag may not be combined with any of the other flags that configure AOF and RDB operations. </w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(12).Range.InsertXML($flagBodyXml)

$p12 = $d.Paragraphs(12)
$p12.Range.InsertParagraphAfter()
$flagSyntaxXml = "<w:p $wNs>" +
  "<w:pPr><w:ind w:firstLine='720'/></w:pPr>" +
  "<w:proofErr w:type='gramStart'/><w:r><w:t>persistence-available</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:t xml:space='preserve'> [(yes)|no]</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(13).Range.InsertXML($flagSyntaxXml)

Write-Output "done"
